# Add Spanish ("label:es") translations to the survey + choices sheets and
# introduce a new "status" (alive/dead) single-select question to the
# plant_form XLSForm, per the authoring change "add 'alive/dead' status to
# form".

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------
# 1. SURVEY sheet: insert a new "label:es" column before the existing
#    "relevant" column (D), shifting "relevant" to column E, then fill in
#    the Spanish labels.
# ---------------------------------------------------------------------
$survey.Columns("D").Insert()

$survey.Range("D1").Value = "label:es"

# username / deviceid "dummy" rows keep the same non-breaking-space note,
# now living in the (shifted) "relevant" column E - nothing else to do
# there since the Insert already moved that content across. D3/D4 stay
# blank (just inherit the row's yellow highlight from the Insert).

# Plant ID (barcode + acc_no_typed)
$survey.Range("D6").Value = "código"
$survey.Range("D7").Value = "código"

# Species
$survey.Range("D9").Value = "especie"

# Location
$survey.Range("D10").Value = "ubicación"

# photo_repeat / photo
$survey.Range("D12").Value = "otra foto"
$survey.Range("D13").Value = "foto"

# ---------------------------------------------------------------------
# 2. SURVEY sheet: insert two new rows above the old row 16 ("note" row),
#    pushing it down to row 18, then add the new "select_one status"
#    question in row 16 (row 15 and 17 stay blank, matching the form's
#    usual blank-row spacing), and repopulate the "note" row (now 18)
#    with its Spanish label.
# ---------------------------------------------------------------------
$survey.Rows("15:16").Insert()

$survey.Range("C14").HorizontalAlignment = -4130  # xlJustify
$survey.Range("D14").HorizontalAlignment = -4130
$survey.Range("C15").HorizontalAlignment = -4130
$survey.Range("D15").HorizontalAlignment = -4130

$survey.Range("A16").Value = "select_one status"
$survey.Range("B16").Value = "alive"
$survey.Range("C16").Value = "Alive"
$survey.Range("D16").Value = "viva"
$survey.Range("C16").HorizontalAlignment = -4130
$survey.Range("D16").HorizontalAlignment = -4130

# The "note" row, now shifted down to row 18.
$survey.Range("D18").Value = "notas"

# ---------------------------------------------------------------------
# 3. CHOICES sheet: add the "label:es" header and a new "status" list
#    with "alive" (TRUE) / "dead" (FALSE) options and Spanish labels.
# ---------------------------------------------------------------------
$choices.Range("D1").Value = "label:es"
# Match the bold header formatting used by the other header cells.
$choices.Range("C1").Copy()
$choices.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

$choices.Range("A2").Value = "status"
$choices.Range("B2").Value = $true
$choices.Range("B2").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$choices.Range("C2").Value = "alive"
$choices.Range("D2").Value = "viva"

$choices.Range("A3").Value = "status"
$choices.Range("B3").Value = $false
$choices.Range("B3").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$choices.Range("C3").Value = "dead"
$choices.Range("D3").Value = "muerta"

# ---------------------------------------------------------------------
# 4. Selection / active-sheet bookkeeping to match the edited workbook's
#    UI state: "choices" becomes the active tab, with A1 selected there,
#    and D18 selected back on "survey" (settings no longer the active tab).
# ---------------------------------------------------------------------
$survey.Range("D18").Select()

$choices.Activate()
$choices.Range("A1").Select()
